$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.0044118118167
$ws.Range("C2").Value = -1.663765130792966
$ws.Range("D2").Value = 0.205149320636871
$ws.Range("E2").Value = 4.754992015228786
$ws.Range("F2").Value = ("1.737961192307855e-05" -as [double])
$ws.Range("G2").Value = 0.7282973682284216
$ws.Range("H2").Value = 0.380665828052104
$ws.Range("I2").Value = 5.851378585728954
$ws.Range("J2").Value = ("1.217454321637733e-14" -as [double])
$ws.Range("K2").Value = 41
$ws.Range("B3").Value = 1.104180648594092
$ws.Range("C3").Value = -1.643308062591759
$ws.Range("D3").Value = 0.1282127434147739
$ws.Range("E3").Value = 2.765832160685457
$ws.Range("F3").Value = ("2.336198482367286e-09" -as [double])
$ws.Range("G3").Value = 0.5571849410297567
$ws.Range("H3").Value = 0.7259421860880489
$ws.Range("I3").Value = 3.377526120081327
$ws.Range("J3").Value = ("-2.723747153747051e-15" -as [double])
$ws.Range("K3").Value = 30
$ws.Range("B4").Value = 1.763124532155738
$ws.Range("C4").Value = -17.44113588547562
$ws.Range("D4").Value = 0.2734058655372659
$ws.Range("E4").Value = 6.31371674319395
$ws.Range("F4").Value = ("1.554158763801481e-07" -as [double])
$ws.Range("G4").Value = 0.008882909693390686
$ws.Range("H4").Value = 0.529180045017565
$ws.Range("I4").Value = 4.85739455003578
$ws.Range("J4").Value = ("-7.925284360401118e-15" -as [double])
$ws.Range("K4").Value = 39
$ws.Range("B5").Value = 1.486936761394673
$ws.Range("C5").Value = -12.18888249123439
$ws.Range("D5").Value = 0.2178369606826689
$ws.Range("E5").Value = 5.167976375962158
$ws.Range("F5").Value = ("2.045380335021619e-07" -as [double])
$ws.Range("G5").Value = 0.02556061569622243
$ws.Range("H5").Value = 0.6246302966024426
$ws.Range("I5").Value = 4.783821989967349
$ws.Range("J5").Value = ("-2.137549396744968e-14" -as [double])
$ws.Range("K5").Value = 30
$ws.Range("B6").Value = 1.708275139345329
$ws.Range("C6").Value = -15.38975200051759
$ws.Range("D6").Value = 0.24521105710697
$ws.Range("E6").Value = 5.578266690555032
$ws.Range("F6").Value = ("1.418164207188524e-07" -as [double])
$ws.Range("G6").Value = 0.01010560478203557
$ws.Range("H6").Value = 0.6341438957377077
$ws.Range("I6").Value = 4.970147232001758
$ws.Range("J6").Value = ("-1.1842378929335e-14" -as [double])
$ws.Range("K6").Value = 30
$ws.Range("B7").Value = 1.498663714003373
$ws.Range("C7").Value = -10.72278843508172
$ws.Range("D7").Value = 0.2148842711123655
$ws.Range("E7").Value = 5.045875125624581
$ws.Range("F7").Value = ("4.123936309411855e-08" -as [double])
$ws.Range("G7").Value = 0.04071624825824748
$ws.Range("H7").Value = 0.5815430501635944
$ws.Range("I7").Value = 5.675229814870046
$ws.Range("J7").Value = ("5.473099451125096e-15" -as [double])
$ws.Range("K7").Value = 37
$ws.Range("B8").Value = 1.071060677624773
$ws.Range("C8").Value = -0.4346948847588554
$ws.Range("D8").Value = 0.210545440308236
$ws.Range("E8").Value = 4.91132860568042
$ws.Range("F8").Value = ("1.673343332618078e-05" -as [double])
$ws.Range("G8").Value = 0.9300417606036931
$ws.Range("H8").Value = 0.4549770928067197
$ws.Range("I8").Value = 5.59253908066523
$ws.Range("J8").Value = ("7.912862284601115e-15" -as [double])
$ws.Range("K8").Value = 33
$ws.Range("B9").Value = 1.66175486910557
$ws.Range("C9").Value = -12.77252199858115
$ws.Range("D9").Value = 0.1151715475268844
$ws.Range("E9").Value = 2.461028939341394
$ws.Range("F9").Value = ("9.108925829895141e-15" -as [double])
$ws.Range("G9").Value = ("1.4962272989251e-05" -as [double])
$ws.Range("H9").Value = 0.8777311088258974
$ws.Range("I9").Value = 2.863486121001196
$ws.Range("J9").Value = ("-1.209068687462751e-14" -as [double])
$ws.Range("K9").Value = 31
